$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new log rows (104, 105) to the feed_logs sheet, following the
# same layout as the existing rows: run_id, rss_url_id, date, response, item_count

$ws.Range("A104").Value = 103
$ws.Range("B104").Value = 1
$ws.Range("C104").Value = "2024-06-17 03:15:33"
$ws.Range("D104").Value = 200
$ws.Range("E104").Value = 11

$ws.Range("A105").Value = 104
$ws.Range("B105").Value = 2
$ws.Range("C105").Value = "2024-06-17 03:15:33"
$ws.Range("D105").Value = 200
$ws.Range("E105").Value = 1
